$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the tooltip text for the "Deep Learning" and "Front-end" rows
$ws.Range("E7").Value = "Keras, TensorFlow"
$ws.Range("E8").Value = "Plotly Dash, HTML, CSS, JavaScript"

# Move the active selection to E28
$ws.Range("E28").Select()
